# Update "想去人数" (interested-count) values in column F across sheets
# 展览 (sheet "展览"), 本地生活 (sheet "本地生活"), 全部类型 (sheet "全部类型")

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 681
$ws1.Range("F3").Value = 48
$ws1.Range("F4").Value = 1990
$ws1.Range("F5").Value = 5815
$ws1.Range("F6").Value = 1632
$ws1.Range("F7").Value = 173
$ws1.Range("F8").Value = 3285
$ws1.Range("F10").Value = 47
$ws1.Range("F11").Value = 1378
$ws1.Range("F12").Value = 4576
$ws1.Range("F13").Value = 1094
$ws1.Range("F14").Value = 1724
$ws1.Range("F15").Value = 2605
$ws1.Range("F17").Value = 51
$ws1.Range("F18").Value = 55
$ws1.Range("F19").Value = 184
$ws1.Range("F21").Value = 1033
$ws1.Range("F24").Value = 17
$ws1.Range("F25").Value = 91
$ws1.Range("F26").Value = 3
$ws1.Range("F28").Value = 10
$ws1.Range("F29").Value = 1127
$ws1.Range("F30").Value = 415
$ws1.Range("F31").Value = 91
$ws1.Range("F32").Value = 208
$ws1.Range("F33").Value = 407
$ws1.Range("F36").Value = 1759
$ws1.Range("F37").Value = 2259
$ws1.Range("F40").Value = 279
$ws1.Range("F41").Value = 640
$ws1.Range("F42").Value = 385
$ws1.Range("F44").Value = 678
$ws1.Range("F45").Value = 34
$ws1.Range("F46").Value = 449
$ws1.Range("F47").Value = 420
$ws1.Range("F48").Value = 233

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 785

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 785
$ws4.Range("F3").Value = 681
$ws4.Range("F4").Value = 48
$ws4.Range("F5").Value = 1990
$ws4.Range("F6").Value = 5815
$ws4.Range("F7").Value = 1632
$ws4.Range("F8").Value = 173
$ws4.Range("F9").Value = 3285
$ws4.Range("F10").Value = 47
$ws4.Range("F11").Value = 1378
$ws4.Range("F12").Value = 4576
$ws4.Range("F13").Value = 1094
$ws4.Range("F14").Value = 1724
$ws4.Range("F17").Value = 51
$ws4.Range("F19").Value = 55
$ws4.Range("F20").Value = 184
$ws4.Range("F23").Value = 1033
$ws4.Range("F26").Value = 91
$ws4.Range("F29").Value = 1127
$ws4.Range("F30").Value = 415
$ws4.Range("F31").Value = 91
$ws4.Range("F32").Value = 208
$ws4.Range("F34").Value = 1759
$ws4.Range("F35").Value = 2259
$ws4.Range("F40").Value = 279
$ws4.Range("F41").Value = 640
$ws4.Range("F42").Value = 385
$ws4.Range("F43").Value = 678
$ws4.Range("F44").Value = 449
$ws4.Range("F45").Value = 420
$ws4.Range("F46").Value = 233
